$d = $word.ActiveDocument

# 1. Remove leading "그리고 " and append "의미가 많은" before "차입니다."
$d.Content.Find.Execute(
    "그리고 오래 전부터 숙박 업체에서 제공되어 왔으며 친한 친구에게 대접하거나 휴식이 필요할 때 마시는 등 역사/문화적으로도 중요한 차입니다.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "오래 전부터 숙박 업체에서 제공되어 왔으며 친한 친구에게 대접하거나 휴식이 필요할 때 마시는 등 역사/문화적으로도 중요한 의미가 많은 차입니다.",
    2
)

# 2. "제품 설명" heading -> "제품 설명:" (only the standalone Heading 1 paragraph whose
#    full text is exactly "제품 설명" -- other occurrences of the same substring
#    elsewhere in the document, and the one already followed by ":" in the table,
#    must stay untouched).
$paraCount = $d.Paragraphs.Count
for ($i = 1; $i -le $paraCount; $i++) {
    $p = $d.Paragraphs($i)
    $trimmed = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($trimmed -eq "제품 설명" -and $p.Style.NameLocal -eq "Heading 1") {
        $p.Range.Find.Execute(
            "제품 설명",
            $true, $false, $false, $false, $false, $true, 1, $false,
            "제품 설명:",
            2
        )
    }
}

# 3. "지속 가능한 공급" -> "지속 가능한 소싱", and "재료를 공급함으로써" -> "재료를 공급받음으로써"
$d.Content.Find.Execute(
    "지속 가능한 공급: 저희는 지속 가능성을 위해 최선을 다하고 있으며, 유기농 농업을 실천하는 소규모 농장에서 재료를 공급함으로써 최고의 품질뿐만 아니라 지구의 건강에도 이바지하고 있습니다.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "지속 가능한 소싱: 저희는 지속 가능성을 위해 최선을 다하고 있으며, 유기농 농업을 실천하는 소규모 농장에서 재료를 공급받음으로써 최고의 품질뿐만 아니라 지구의 건강에도 이바지하고 있습니다.",
    2
)

# 4. "단 음료" -> "단 간식"
$d.Content.Find.Execute(
    "또한 차이 티는 단 음료를 함께 즐기기를 좋아하는 라틴 아메리카 소비자의 생활 방식과 취향에도 적합합니다.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "또한 차이 티는 단 간식을 함께 즐기기를 좋아하는 라틴 아메리카 소비자의 생활 방식과 취향에도 적합합니다.",
    2
)
